# Add a new contract-note trade as the most recent entry in "Trading History".
# This pushes the existing rows 5-11 down to 6-12 and inserts a brand-new
# BSE "Buy" row at row 5 (date 2026-02-10 / serial 46063) that also carries
# STT (H) and ADD CHRG (I) values - data not present on the older rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trading History")

# Insert a blank row above row 5; existing rows 5-11 shift down to 6-12.
$ws.Rows.Item(5).Insert()

# Excel's Insert() copies the formatting of the row above (the bold header
# row) onto the newly inserted row. Wipe that back to the plain look used
# by the rest of the data rows before writing the new values.
$ws.Range("A5:N5").Clear()
$ws.Range("W5:AB5").Clear()

# DATE - keep the same date number format used by the other rows in col A.
$ws.Cells.Item(5, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(5, 1).Value = 46063

$ws.Cells.Item(5, 2).Value = "BSE"                 # EXCH
$ws.Cells.Item(5, 3).Value = "Buy"                 # ACTION
$ws.Cells.Item(5, 4).Value = 1                      # QTY
$ws.Cells.Item(5, 5).Value = 11645                  # PRICE
$ws.Cells.Item(5, 6).Value = 11727.59               # COST
$ws.Cells.Item(5, 7).Value = "CN#252611730667"      # REMARKS
$ws.Cells.Item(5, 8).Value = 11.62                  # STT
$ws.Cells.Item(5, 9).Value = 70.97                  # ADD CHRG
$ws.Cells.Item(5, 10).Formula = "=Index!`$C`$2"     # Current Price
